$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1156.5333
$ws.Range("I6").Value = 889.1429000000001
$ws.Range("K6").Value = 2667.4287
$ws.Range("M6").Value = -2555.4287

$ws.Range("H8").Value = 339.77585
$ws.Range("I8").Value = 499.4
$ws.Range("K8").Value = 1498.2
$ws.Range("M8").Value = -1359.2

$ws.Range("H12").Value = 209
$ws.Range("I12").Value = 215
$ws.Range("J12").Value = 149
$ws.Range("K12").Value = 215
$ws.Range("L12").Value = 149
$ws.Range("M12").Value = -45
$ws.Range("N12").Value = -489

$ws.Range("H17").Value = 1593.2727
$ws.Range("J17").Value = 1652.6
$ws.Range("L17").Value = 4957.799999999999
$ws.Range("N17").Value = -5293.799999999999

$ws.Range("H99").Value = 1585
$ws.Range("I99").Value = 309.5
$ws.Range("J99").Value = 2070.9048
$ws.Range("K99").Value = 928.5
$ws.Range("L99").Value = 6212.714399999999
$ws.Range("M99").Value = 569.5
$ws.Range("N99").Value = -9208.714399999999

$ws.Range("H100").Value = 5386.9165
$ws.Range("J100").Value = 1651.75
$ws.Range("L100").Value = 1651.75
$ws.Range("N100").Value = -2733.75

$ws.Range("H111").Value = 7833.3335
$ws.Range("I111").Value = 6750
$ws.Range("J111").Value = 10000
$ws.Range("K111").Value = 20250
$ws.Range("L111").Value = 30000
$ws.Range("M111").Value = -17183
$ws.Range("N111").Value = -36134

$ws.Range("H116").Value = 15320.725
$ws.Range("I116").Value = 18257.666
$ws.Range("K116").Value = 18257.666
$ws.Range("M116").Value = -14815.666

$ws.Range("H118").Value = 910.6
$ws.Range("I118").Value = 547.0714
$ws.Range("K118").Value = 1641.2142
$ws.Range("M118").Value = 15.78579999999988

$ws.Range("H127").Value = 1660.5714
$ws.Range("J127").Value = 3275
$ws.Range("L127").Value = 9825
$ws.Range("N127").Value = -19745

$ws.Range("H129").Value = 1850.9166
$ws.Range("I129").Value = 1357.4445
$ws.Range("K129").Value = 4072.3335
$ws.Range("M129").Value = 927.6664999999998

$ws.Range("H137").Value = 25934.166
$ws.Range("I137").Value = 60779.2
$ws.Range("J137").Value = 1044.8572
$ws.Range("K137").Value = 182337.6
$ws.Range("L137").Value = 3134.5716
$ws.Range("M137").Value = -179787.6
$ws.Range("N137").Value = -8234.571599999999

$ws.Range("H138").Value = 32857.758
$ws.Range("I138").Value = 2172.35
$ws.Range("K138").Value = 6517.049999999999
$ws.Range("M138").Value = -1377.049999999999

$ws.Range("H141").Value = 2978.5
$ws.Range("I141").Value = 2978.5
$ws.Range("K141").Value = 8935.5
$ws.Range("M141").Value = -3755.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3503.8667
$ws.Range("I5").Value = 4724.909
$ws.Range("J5").Value = 146
$ws.Range("K5").Value = 4724.909
$ws.Range("L5").Value = 146
$ws.Range("M5").Value = -4612.909
$ws.Range("N5").Value = -370

$ws.Range("H110").Value = 7143.3335
$ws.Range("I110").Value = 7143.3335
$ws.Range("K110").Value = 7143.3335
$ws.Range("M110").Value = -5098.3335

$ws.Range("H132").Value = 1239.0566
$ws.Range("I132").Value = 1027.8478
$ws.Range("J132").Value = 2627
$ws.Range("K132").Value = 3083.5434
$ws.Range("L132").Value = 7881
$ws.Range("M132").Value = -553.5434
$ws.Range("N132").Value = -12941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3503.8667
$ws.Range("I4").Value = 4724.909
$ws.Range("J4").Value = 146
$ws.Range("K4").Value = 4724.909
$ws.Range("L4").Value = 146
$ws.Range("M4").Value = -4609.909
$ws.Range("N4").Value = -376

$ws.Range("H64").Value = 1323.6364
$ws.Range("I64").Value = 1743.3334
$ws.Range("J64").Value = 1166.25
$ws.Range("K64").Value = 1743.3334
$ws.Range("L64").Value = 1166.25
$ws.Range("M64").Value = -1518.3334
$ws.Range("N64").Value = -1616.25

$ws.Range("H67").Value = 1323.6364
$ws.Range("I67").Value = 1743.3334
$ws.Range("J67").Value = 1166.25
$ws.Range("K67").Value = 1743.3334
$ws.Range("L67").Value = 1166.25
$ws.Range("M67").Value = -963.3334
$ws.Range("N67").Value = -2726.25

$ws.Range("H105").Value = 1981.9615
$ws.Range("I105").Value = 1312.6111
$ws.Range("K105").Value = 1312.6111
$ws.Range("M105").Value = 434.3888999999999

$ws.Range("H134").Value = 2445.8845
$ws.Range("I134").Value = 2031.8572
$ws.Range("K134").Value = 6095.571599999999
$ws.Range("M134").Value = -3560.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6476.2856
$ws.Range("I62").Value = 5683.5
$ws.Range("J62").Value = 7533.3335
$ws.Range("K62").Value = 5683.5
$ws.Range("L62").Value = 7533.3335
$ws.Range("M62").Value = -5059.5
$ws.Range("N62").Value = -8781.333500000001

$ws.Range("H65").Value = 6476.2856
$ws.Range("I65").Value = 5683.5
$ws.Range("J65").Value = 7533.3335
$ws.Range("K65").Value = 28417.5
$ws.Range("L65").Value = 37666.6675
$ws.Range("M65").Value = -25297.5
$ws.Range("N65").Value = -43906.6675

$ws.Range("H94").Value = 3369.6
$ws.Range("J94").Value = 4018.8
$ws.Range("L94").Value = 4018.8
$ws.Range("N94").Value = -4920.8

$ws.Range("H134").Value = 3628.5217
$ws.Range("I134").Value = 3029.4375
$ws.Range("K134").Value = 9088.3125
$ws.Range("M134").Value = -6553.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 1600
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 9600
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -9487
$ws.Range("N2").Value = -1426

$ws.Range("H38").Value = 90909130
$ws.Range("I38").Value = 63.666668
$ws.Range("J38").Value = 125000024
$ws.Range("K38").Value = 191.000004
$ws.Range("L38").Value = 375000072
$ws.Range("M38").Value = 155.999996
$ws.Range("N38").Value = -375000766

$ws.Range("H129").Value = 3476.2856
$ws.Range("I129").Value = 2827.1667
$ws.Range("J129").Value = 3963.125
$ws.Range("K129").Value = 8481.500100000001
$ws.Range("L129").Value = 11889.375
$ws.Range("M129").Value = -3481.500100000001
$ws.Range("N129").Value = -21889.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136

$ws.Range("H113").Value = 1892.3334
$ws.Range("I113").Value = 1881.8
$ws.Range("J113").Value = 1945
$ws.Range("K113").Value = 1881.8
$ws.Range("L113").Value = 1945
$ws.Range("M113").Value = 288.2
$ws.Range("N113").Value = -6285

$ws.Range("H132").Value = 2018.6052
$ws.Range("I132").Value = 1821.6818
$ws.Range("J132").Value = 2289.375
$ws.Range("K132").Value = 5465.0454
$ws.Range("L132").Value = 6868.125
$ws.Range("M132").Value = -2935.0454
$ws.Range("N132").Value = -11928.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 150000
$ws.Range("J81").Value = 150000
$ws.Range("L81").Value = 150000
$ws.Range("N81").Value = -151996

$ws.Range("H84").Value = 150000
$ws.Range("J84").Value = 150000
$ws.Range("L84").Value = 450000
$ws.Range("N84").Value = -459984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3400.4285
$ws.Range("I126").Value = 2817.7273
$ws.Range("K126").Value = 8453.1819
$ws.Range("M126").Value = -5983.1819

$ws.Range("H136").Value = 22101.719
$ws.Range("I136").Value = 24473.393
$ws.Range("K136").Value = 73420.179
$ws.Range("M136").Value = -70870.179
